# Updates the cryptos list (Coin/Link/Price/Volume(1h)) to the latest
# scraped snapshot, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''68.503.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.41%  '
$ws.Range("D3").Value = '''2.456.87'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.25%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''559.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.69%  '
$ws.Range("D6").Value = '''163.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -1.65%  '
$ws.Range("D9").Value = '''2.456.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.26%  '
$ws.Range("D10").Value = '''0.149'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.65%  '
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("D12").Value = '''0.336'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.97%  '
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("D14").Value = '''2.910.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.26%  '
$ws.Range("D15").Value = '''68.399.14'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("D16").Value = '''0.0000168'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.20%  '
$ws.Range("D17").Value = '''23.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.54%  '
$ws.Range("D18").Value = '''2.457.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.79%  '
$ws.Range("D19").Value = '''10.91'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.95%  '
$ws.Range("D20").Value = '''7.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.74%  '
$ws.Range("D21").Value = '''340.88'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.36%  '
$ws.Range("D22").Value = '''3.76'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.72%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = '''1.85'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.94%  '
$ws.Range("D25").Value = '''67.50'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.05%  '
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = '''3.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.97%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '''2.584.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.46%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '''1.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.86%  '
$ws.Range("D29").Value = '''8.04'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.44%  '
$ws.Range("D30").Value = '''0.0₃0826'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.46%  '
$ws.Range("D31").Value = '''7.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.55%  '
$ws.Range("D32").Value = '''1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.07%  '
$ws.Range("D33").Value = '''428.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.25%  '
$ws.Range("E34").Value = '  -2.73%  '
$ws.Range("D35").Value = '''1.65'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.29%  '
$ws.Range("D36").Value = '''157.21'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.80%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").Value = '''0.109'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.76%  '
$ws.Range("D40").Value = '''17.83'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.09%  '
$ws.Range("D41").Value = '''0.304'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.05%  '
$ws.Range("D42").Value = '''4.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.06%  '
$ws.Range("D43").Value = '''1.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.57%  '
$ws.Range("E44").Value = '  +0.63%  '
$ws.Range("D45").Value = '''2.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.78%  '
$ws.Range("D46").Value = '''134.40'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.72%  '
$ws.Range("D47").Value = '''3.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.23%  '
$ws.Range("E48").Value = '  -2.05%  '
$ws.Range("D49").Value = '''0.479'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.32%  '
$ws.Range("D50").Value = '''0.562'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.83%  '
$ws.Range("E51").Value = '  -1.48%  '
